# "typeahead and description enlarged functionality"
# Apply the edits described by the diff:
#  - B9 activity description text changed
#  - E9,F9,G9 (and the mirrored Total Hours row E18,F18,G18) hours changed from "1" to "11"
#  - D19 total-hours-for-the-week changed from "7" to "37"
#  - B22 project manager name changed
#  - B25 client manager name changed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Activity / status-call description text (B9)
$ws.Range("B9").Value = "q ( bvcvbcbbc ) "

# 2) Hours entered Mon-Wed changed from "1" to "11" on both the activity row (9)
#    and the mirrored Total Hours row (18). These are stored as text in the
#    workbook, so force a text number format before writing the values so the
#    numeric-looking strings aren't silently re-typed as numbers. (Multi-area
#    ranges only apply NumberFormat to the first area, so set each
#    contiguous range separately.)
$ws.Range("E9:G9").NumberFormat = "@"
$ws.Range("E18:G18").NumberFormat = "@"

$ws.Range("E9").Value = "11"
$ws.Range("F9").Value = "11"
$ws.Range("G9").Value = "11"
$ws.Range("E18").Value = "11"
$ws.Range("F18").Value = "11"
$ws.Range("G18").Value = "11"

# 3) Total hours for the week (D19): 7 -> 37
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37"

# 4) Offshore consultant's project manager name (B22)
$ws.Range("B22").Value = "Offshore Consultant's Project Manager's Name :Tanuj Khaturia"

# 5) Client manager name (B25)
$ws.Range("B25").Value = "Omar Colon"
